# Fill in the first two replies (DQ1 response 1 and DQ2 response 1) for week 2
# with the standard 15-minute (1.0416666666666666E-2 day) actual completion time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week2")

$quarterHour = 0.010416666666666666

# C8  -> "DQ1 response 1" actual time length to complete
$ws.Range("C8").Value = $quarterHour

# C13 -> "DQ2 response 1" actual time length to complete
$ws.Range("C13").Value = $quarterHour

# Move the active selection to C14, matching the author's final cursor position
$ws.Activate()
$ws.Range("C14").Select()
